$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1126.5
$ws.Range("I2").Value = 178.42857
$ws.Range("J2").Value = 1568.9333
$ws.Range("K2").Value = 178.42857
$ws.Range("L2").Value = 1568.9333
$ws.Range("M2").Value = -65.42857000000001
$ws.Range("N2").Value = -1794.9333
$ws.Range("H9").Value = 256.42105
$ws.Range("I9").Value = 247.81818
$ws.Range("K9").Value = 247.81818
$ws.Range("M9").Value = -78.81818000000001
$ws.Range("H15").Value = 2837.2354
$ws.Range("I15").Value = 2837.2354
$ws.Range("K15").Value = 8511.706200000001
$ws.Range("M15").Value = -8342.706200000001
$ws.Range("H43").Value = 4140.778
$ws.Range("J43").Value = 4402.533
$ws.Range("L43").Value = 4402.533
$ws.Range("N43").Value = -4540.533
$ws.Range("H51").Value = 23200
$ws.Range("I51").Value = 23200
$ws.Range("K51").Value = 23200
$ws.Range("M51").Value = -22716
$ws.Range("H70").Value = 1615.4
$ws.Range("J70").Value = 1694.6666
$ws.Range("L70").Value = 5083.9998
$ws.Range("N70").Value = -5623.9998
$ws.Range("H73").Value = 1615.4
$ws.Range("J73").Value = 1694.6666
$ws.Range("L73").Value = 5083.9998
$ws.Range("N73").Value = -6955.9998
$ws.Range("H74").Value = 13853.0625
$ws.Range("I74").Value = 7240.8184
$ws.Range("K74").Value = 7240.8184
$ws.Range("M74").Value = -6304.8184
$ws.Range("H77").Value = 13853.0625
$ws.Range("I77").Value = 7240.8184
$ws.Range("K77").Value = 36204.092
$ws.Range("M77").Value = -31524.092
$ws.Range("H98").Value = 1302.3784
$ws.Range("I98").Value = 1265.931
$ws.Range("K98").Value = 1265.931
$ws.Range("M98").Value = 232.069
$ws.Range("H112").Value = 3180.5454
$ws.Range("I112").Value = 1999
$ws.Range("J112").Value = 3298.7
$ws.Range("K112").Value = 5997
$ws.Range("L112").Value = 9896.099999999999
$ws.Range("M112").Value = -4889
$ws.Range("N112").Value = -12112.1
$ws.Range("H122").Value = 1302.3784
$ws.Range("I122").Value = 1265.931
$ws.Range("K122").Value = 3797.793
$ws.Range("M122").Value = -1347.793
$ws.Range("H138").Value = 3428.1785
$ws.Range("I138").Value = 3470.1853
$ws.Range("J138").Value = 2294
$ws.Range("K138").Value = 10410.5559
$ws.Range("L138").Value = 6882
$ws.Range("M138").Value = -5270.555899999999
$ws.Range("N138").Value = -17162

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 32407.344
$ws.Range("I32").Value = 34531.227
$ws.Range("J32").Value = 15947.25
$ws.Range("K32").Value = 34531.227
$ws.Range("L32").Value = 15947.25
$ws.Range("M32").Value = -34244.227
$ws.Range("N32").Value = -16521.25
$ws.Range("H45").Value = 2289.682
$ws.Range("I45").Value = 1270.6923
$ws.Range("J45").Value = 3761.5557
$ws.Range("K45").Value = 1270.6923
$ws.Range("L45").Value = 3761.5557
$ws.Range("M45").Value = -893.6922999999999
$ws.Range("N45").Value = -4515.5557
$ws.Range("H62").Value = 30249
$ws.Range("J62").Value = 30249
$ws.Range("L62").Value = 30249
$ws.Range("N62").Value = -31497
$ws.Range("H65").Value = 30249
$ws.Range("J65").Value = 30249
$ws.Range("L65").Value = 90747
$ws.Range("N65").Value = -96987
$ws.Range("H74").Value = 1514.8334
$ws.Range("I74").Value = 1091.25
$ws.Range("K74").Value = 1091.25
$ws.Range("M74").Value = -217.25
$ws.Range("H77").Value = 1514.8334
$ws.Range("I77").Value = 1091.25
$ws.Range("K77").Value = 5456.25
$ws.Range("M77").Value = -1088.25
$ws.Range("H110").Value = 3866.25
$ws.Range("I110").Value = 3488.3333
$ws.Range("K110").Value = 3488.3333
$ws.Range("M110").Value = -1443.3333
$ws.Range("H122").Value = 2052.476
$ws.Range("I122").Value = 1858.2941
$ws.Range("J122").Value = 2877.75
$ws.Range("K122").Value = 5574.8823
$ws.Range("L122").Value = 8633.25
$ws.Range("M122").Value = -3124.8823
$ws.Range("N122").Value = -13533.25
$ws.Range("H128").Value = 28999
$ws.Range("J128").Value = 28999
$ws.Range("L128").Value = 28999
$ws.Range("N128").Value = -38959

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H39").Value = 13350.667
$ws.Range("J39").Value = 13350.667
$ws.Range("L39").Value = 13350.667
$ws.Range("N39").Value = -14128.667
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").Value = $null
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = $null
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").Value = $null
$ws.Range("H80").Value = 227.25
$ws.Range("J80").Value = 220.63637
$ws.Range("L80").Value = 220.63637
$ws.Range("N80").Value = -2216.63637
$ws.Range("H83").Value = 227.25
$ws.Range("J83").Value = 220.63637
$ws.Range("L83").Value = 1103.18185
$ws.Range("N83").Value = -11087.18185
$ws.Range("H86").Value = 4388.857
$ws.Range("I86").Value = 3499.375
$ws.Range("K86").Value = 3499.375
$ws.Range("M86").Value = -2376.375
$ws.Range("H89").Value = 4388.857
$ws.Range("I89").Value = 3499.375
$ws.Range("K89").Value = 17496.875
$ws.Range("M89").Value = -11880.875
$ws.Range("H105").Value = 3672.3635
$ws.Range("I105").Value = 3488.5
$ws.Range("J105").Value = 4499.75
$ws.Range("K105").Value = 3488.5
$ws.Range("L105").Value = 4499.75
$ws.Range("M105").Value = -1741.5
$ws.Range("N105").Value = -7993.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9301.23
$ws.Range("I31").Value = 8537.166999999999
$ws.Range("J31").Value = 9956.143
$ws.Range("K31").Value = 8537.166999999999
$ws.Range("L31").Value = 9956.143
$ws.Range("M31").Value = -8242.166999999999
$ws.Range("N31").Value = -10546.143
$ws.Range("H34").Value = 9301.23
$ws.Range("I34").Value = 8537.166999999999
$ws.Range("J34").Value = 9956.143
$ws.Range("K34").Value = 8537.166999999999
$ws.Range("L34").Value = 9956.143
$ws.Range("M34").Value = -8335.166999999999
$ws.Range("N34").Value = -10360.143
$ws.Range("H74").Value = 25993.334
$ws.Range("J74").Value = 27192
$ws.Range("L74").Value = 27192
$ws.Range("N74").Value = -28940
$ws.Range("H77").Value = 25993.334
$ws.Range("J77").Value = 27192
$ws.Range("L77").Value = 81576
$ws.Range("N77").Value = -90312
$ws.Range("H86").Value = 9894.5
$ws.Range("I86").Value = 9894.5
$ws.Range("K86").Value = 9894.5
$ws.Range("M86").Value = -8771.5
$ws.Range("H89").Value = 9894.5
$ws.Range("I89").Value = 9894.5
$ws.Range("K89").Value = 49472.5
$ws.Range("M89").Value = -43856.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 630.75
$ws.Range("I122").Value = 393
$ws.Range("J122").Value = 868.5
$ws.Range("K122").Value = 3537
$ws.Range("L122").Value = 7816.5
$ws.Range("M122").Value = -1087
$ws.Range("N122").Value = -12716.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4914.8184
$ws.Range("I70").Value = 4383.125
$ws.Range("K70").Value = 4383.125
$ws.Range("M70").Value = -4113.125
$ws.Range("H73").Value = 4914.8184
$ws.Range("I73").Value = 4383.125
$ws.Range("K73").Value = 4383.125
$ws.Range("M73").Value = -3447.125
$ws.Range("H102").Value = 3296.44
$ws.Range("I102").Value = 2325.9375
$ws.Range("J102").Value = 5021.778
$ws.Range("K102").Value = 2325.9375
$ws.Range("L102").Value = 5021.778
$ws.Range("M102").Value = -703.9375
$ws.Range("N102").Value = -8265.778
$ws.Range("H122").Value = 2898.6453
$ws.Range("I122").Value = 2187
$ws.Range("J122").Value = 5338.5713
$ws.Range("K122").Value = 6561
$ws.Range("L122").Value = 16015.7139
$ws.Range("M122").Value = -4111
$ws.Range("N122").Value = -20915.7139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2742.7144
$ws.Range("J93").Value = 3149.75
$ws.Range("L93").Value = 3149.75
$ws.Range("N93").Value = -5645.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").Value = $null
$ws.Range("H34").Value = 27999.5
$ws.Range("J34").Value = 45999
$ws.Range("L34").Value = 45999
$ws.Range("N34").Value = -46405
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").Value = $null
$ws.Range("H107").Value = 1376.9445
$ws.Range("I107").Value = 635.5454999999999
$ws.Range("K107").Value = 1906.6365
$ws.Range("M107").Value = 13.36350000000016
$ws.Range("H132").Value = 48793.184
$ws.Range("I132").Value = 48793.184
$ws.Range("K132").Value = 146379.552
$ws.Range("M132").Value = -143849.552
$ws.Range("H136").Value = 2897.8572
$ws.Range("I136").Value = 2589.1667
$ws.Range("K136").Value = 7767.500100000001
$ws.Range("M136").Value = -5217.500100000001
